# Review_272.docx update:
#  - paragraph 1: date 10.08.24 -> 09.08.24
#  - paragraph 2: new title, drop the trailing line break
#  - paragraphs 3-6: replace the body text with the new review text
#  - paragraph 7: new arxiv link
#
# Each paragraph is rewritten by deleting the paragraph's own content
# (its Range minus the trailing paragraph-mark character) and then
# inserting the replacement text into the now-collapsed range. This
# swaps the run text cleanly (and drops stray children such as <w:br/>)
# without disturbing paragraph boundaries or styles.
#
# Paragraphs are processed from last to first so that the character
# offsets of not-yet-edited (earlier) paragraphs stay valid while later
# ones are being resized.

$d = $word.ActiveDocument

function Set-ParagraphText($doc, [int]$index, [string]$newText) {
    $para = $doc.Paragraphs.Item($index)
    $rg = $doc.Range($para.Range.Start, $para.Range.End - 1)
    $rg.Delete()
    if ($para.Range.Start -eq 0) {
        # InsertAfter on a collapsed range sitting at document position 0
        # lands at the end of the story instead of at the start, so the
        # very first paragraph needs InsertBefore instead.
        $rg.InsertBefore($newText)
    } else {
        $rg.InsertAfter($newText)
    }
}

Set-ParagraphText $d 7 "https://arxiv.org/abs/2408.03314"

Set-ParagraphText $d 6 "מה השיטה העדיפה לרמת ביצועים אופטימלית בהינתן תקציב חישוב נתון (FLOps) - זו השאלה שהמאמר מנסה לענות עליה ויש תוצאות מעניינות (לדעתי)"

Set-ParagraphText $d 5 'יש שיטות איטרטיביות אחרות כמו במאמר "Consistency LLMs" שסקרתי לפני כמה שבועות. הוצעו גם שיטות שמשערכות את ״איכות״ התשובה המגונרטת (עם מודל מאומן נוסף) שמאפשר לבחור את התשובה הכי טובה מכמה תשובות מגונרטות (או להפסיק את יצירת התשובה אם רואים שהיא לא ״בכיוון). כל שיטה כזו דורשת משאבי חישוב שונים שתלויים גם בהייפרפרמטרים של השיטה.'

Set-ParagraphText $d 4 "המאמר שואל האם ניתן לנסח חוקי סקיילנג דומים עבור האינפרנס, כלומר מה הביצועים המקסימליים שניתן להפיק בהינתן כמות משאבי חישוב נתונה. הרי יש כמה שיטות לבצע אינפרנס של מודל השפה ויש כמה פרמטרים חשובים של האינפרנס המשפיעים בצורה משמעותית על הביצועים. למשל יש שיטה הנקראת beam search שיוצרת בכל חיזוי של טוקן M סדרות טוקנים בעלי נראות (likelihood) הגבוהה ביותר. קיימות שיטות beam search עם מספר הסדרות השמורות לא קבוע ותלוי במספר הטוקן המגונרט."

Set-ParagraphText $d 3 "בטח שמעתם על חוקי הסקיילינג של מודלי שפה. חוקים אלו מיועדים למציאת ״קונפיגורציה״ אופטימלית לאימון מודלי שפה.  חוקי סקליינג מקשרים ערך של פונקציית לוס (ניתן להגדיר אותו בכמה אופנים) שניתן להשיגו עבור גודל מודל, גודל סט האימון וכמות משאבי החישוב (FLOps) המוקצית לאימון. "

Set-ParagraphText $d 2 "Scaling LLM Test-Time Compute Optimally can be More Effective than Scaling Model Parameters"

Set-ParagraphText $d 1 "⚡️🚀המאמר היומי של מייק 09.08.24: ⚡️🚀"
